# "why not both round 1" - fill in the Round 1 winners that were left
# showing the placeholder "p" value in column O for the bottom half of
# the bracket (rows 19-33), matching the picks already recorded in
# column P for each match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("O19").Value = "Walrus"
$ws.Range("O21").Value = "Echidna"
$ws.Range("O23").Value = "Serval"
$ws.Range("O25").Value = "Pangolin"
$ws.Range("O27").Value = "Therapsid"
$ws.Range("O29").Value = "Spotted Salamander"
$ws.Range("O31").Value = "Hairy Frogfish"
$ws.Range("O33").Value = "Swordfish"

# Move the active selection / scroll position the same way the author
# left it after making the picks.
$ws.Activate()
$ws.Range("O36").Select()
